$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$features = @(
    "concessão",
    "deferi",
    "liminar hc",
    "república concessão",
    "ordem prisão",
    "assuntos COMPETÊNCIA DO MP",
    "opina concessão",
    "senha relatório",
    "sobrestamento",
    "justiça indeferiu",
    "violência grave",
    "enunciado súmula",
    "autoridade judiciária",
    "xliii constituição",
    "outro motivo",
    "Relator_OCTAVIO GALLOTTI",
    "deferido",
    "resumiu apreciou",
    "código senha",
    "apreciou presente",
    "opina deferimento",
    "ministro gilson",
    "restritiva direitos",
    "deferimento",
    "efeitos ordem",
    "medidas cautelares",
    "deduzida",
    "cento",
    "indevidas",
    "arquivado",
    "http sob",
    "concessão ordem",
    "deferida",
    "indeferiu liminarmente",
    "http",
    "liberdade restritiva",
    "relatório http",
    "vedação liberdade",
    "submetida",
    "cautelares previstas",
    "sob código",
    "flagrante preventiva",
    "inconstitucional",
    "stj indeferiu",
    "afastamento enunciado",
    "senha",
    "contornos",
    "liminar assessor",
    "decretada desfavor",
    "gabinete prestou",
    "estado paulo",
    "definitivo após",
    "arquivado definitivo",
    "suspender",
    "senha primeira",
    "decisão proferida",
    "efêmero",
    "assuntos LICITAÇÕES",
    "opinou concessão",
    "liminar suspender",
    "precário",
    "deferi pedido",
    "implicou deferimento",
    "requer medida",
    "análise pedido",
    "suspender efeitos",
    "fiança",
    "ficou assim",
    "precário efêmero",
    "enunciado",
    "informado",
    "assuntos FIANÇA",
    "proferida ministro",
    "assim revelou",
    "resumida",
    "revelou contornos",
    "assim resumida",
    "deferida assessoria",
    "impetração eis",
    "liminar deferida",
    "processo formalizado",
    "espécie ficou",
    "informado análise",
    "eis informado",
    "contornos impetração",
    "assuntos TRANSFERÊNCIA DE PRESO",
    "assuntos INDEFERIMENTO",
    "assuntos DESCLASSIFICAÇÃO",
    "formalizado ato",
    "superação súmula",
    "assuntos FURTO (ART. 155)",
    "alcance",
    "Relator_ILMAR GALVÃO",
    "aurélio decisão",
    "resumida prisão",
    "deserção",
    "fundamentos insubsistência",
    "liminar espécie",
    "preventiva fundamentos",
    "assuntos HABEAS CORPUS - LIBERATÓRIO"
)

$importances = @(
    [double]"0.03943627219794957",
    [double]"0.03882963684287932",
    [double]"0.03662977004559582",
    [double]"0.03150191172787244",
    [double]"0.03079242482761155",
    [double]"0.02888210675098866",
    [double]"0.02789687800920582",
    [double]"0.02732471745915093",
    [double]"0.02719979969737903",
    [double]"0.02686392226343103",
    [double]"0.02562121723432194",
    [double]"0.02560574144406425",
    [double]"0.02378913617963958",
    [double]"0.02377441625243438",
    [double]"0.02376524928008462",
    [double]"0.02034016661823855",
    [double]"0.01953334947005041",
    [double]"0.01931619868240219",
    [double]"0.01921377237481184",
    [double]"0.01890359086036458",
    [double]"0.01847873440707344",
    [double]"0.01752647336230672",
    [double]"0.01685491483248903",
    [double]"0.01671966393214012",
    [double]"0.01635280171439729",
    [double]"0.01607357300009383",
    [double]"0.01577045041769757",
    [double]"0.01568287597373564",
    [double]"0.01563834586482164",
    [double]"0.0150562090334635",
    [double]"0.01451736262860866",
    [double]"0.01365428972386297",
    [double]"0.01333873659585741",
    [double]"0.01328201725087755",
    [double]"0.01327306439680962",
    [double]"0.01319918873628093",
    [double]"0.01312190427081805",
    [double]"0.0126310848858902",
    [double]"0.01262710065610201",
    [double]"0.01254165421097431",
    [double]"0.0125099535383994",
    [double]"0.01092318651250999",
    [double]"0.01089555967036921",
    [double]"0.01061919418467707",
    [double]"0.01047569335347573",
    [double]"0.0102329517617135",
    [double]"0.009450465226043064",
    [double]"0.008078938485714024",
    [double]"0.008005554513555228",
    [double]"0.007936799018765513",
    [double]"0.007427811729627588",
    [double]"0.00533576765185011",
    [double]"0.005124710627229247",
    [double]"0.004152456407958737",
    [double]"0.004053805159999832",
    [double]"0.003923681810480437",
    [double]"0.003504714603073834",
    [double]"0.003311318015153255",
    [double]"0.003185851563499837",
    [double]"0.002761984981374526",
    [double]"0.002732966450261067",
    [double]"0.002414572782649035",
    [double]"0.002382456803106262",
    [double]"0.002321808534278295",
    [double]"0.002304376321535941",
    [double]"0.001958321403605697",
    [double]"0.001944308401443623",
    [double]"0.001564512521324016",
    [double]"0.001137301357017755",
    [double]"0.001114712912154443",
    [double]"0.0008352831567786553",
    [double]"0.0007735003558620728",
    [double]"0.0007281105438923458",
    [double]"0.0007131398757218936",
    [double]"0.000666925408866472",
    [double]"0.000560454167563593",
    [double]"0.0005335322502163179",
    [double]"0.0004712080694227272",
    [double]"0.0003998235032391884",
    [double]"0.0003788202562451232",
    [double]"0.0003191344821202926",
    [double]"0.0002736075144466114",
    [double]"8.738806484581303E-16",
    [double]"5.136380628675219E-16",
    [double]"4.793458267680357E-17",
    [double]"0",
    [double]"0",
    [double]"0",
    [double]"0",
    [double]"0",
    [double]"0",
    [double]"0",
    [double]"0",
    [double]"0",
    [double]"0",
    [double]"0",
    [double]"0",
    [double]"0",
    [double]"0",
    [double]"0"
)

for ($i = 0; $i -lt $features.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $features[$i]
    $ws.Cells.Item($row, 2).Value = $importances[$i]
}
